# Daily attendance processing - 2025-12-06 03:50:53
# Swap the order of the "Recorded By" names in column G so that the
# dnasr281@gmail.com entry is listed first, followed by the other
# recorder (System / admin@admin.com), for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$rows = @(3,6,10,11,12,13,14,15,17,18,19,20,21,22,24,26,29,32,36,37,38,39,40,41,43,44,45,46,47,48,50,52,55,58,62,63,64,65,66,67,69,70,71,72,73,74,76,78,83,84,85,86,87,90,92,93,94,96,99,101,109,110,111,112,113,116,118,119,120,122,125,127,135,136,137,138,139,142,144,145,146,148,151,153)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value2
    if ($val) {
        $parts = $val -split ", "
        if ($parts.Count -eq 2) {
            $newVal = $parts[1] + ", " + $parts[0]
            $cell.Value2 = $newVal
        }
    }
}
